$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 47497.848
$ws.Range("J17").Value = 47497.848
$ws.Range("L17").Value = 142493.544
$ws.Range("N17").Value = -142829.544
$ws.Range("H58").Value = 60320.766
$ws.Range("I58").Value = 311.625
$ws.Range("J58").Value = 113662.22
$ws.Range("K58").Value = 934.875
$ws.Range("L58").Value = 340986.66
$ws.Range("M58").Value = -784.875
$ws.Range("N58").Value = -341286.66
$ws.Range("H82").Value = 2790
$ws.Range("I82").Value = 185
$ws.Range("K82").Value = 555
$ws.Range("M82").Value = -149
$ws.Range("H85").Value = 2790
$ws.Range("I85").Value = 185
$ws.Range("K85").Value = 555
$ws.Range("M85").Value = 849
$ws.Range("H132").Value = 2987365
$ws.Range("I132").Value = 3775843.2
$ws.Range("J132").Value = 2412.2856
$ws.Range("K132").Value = 11327529.6
$ws.Range("L132").Value = 7236.8568
$ws.Range("M132").Value = -11324999.6
$ws.Range("N132").Value = -12296.8568
$ws.Range("H137").Value = 3778.0232
$ws.Range("I137").Value = 4345.069
$ws.Range("J137").Value = 2603.4285
$ws.Range("K137").Value = 13035.207
$ws.Range("L137").Value = 7810.2855
$ws.Range("M137").Value = -10485.207
$ws.Range("N137").Value = -12910.2855
$ws.Range("H141").Value = 445249.72
$ws.Range("I141").Value = 3261.6458
$ws.Range("J141").Value = 2566792.5
$ws.Range("K141").Value = 9784.937399999999
$ws.Range("L141").Value = 7700377.5
$ws.Range("M141").Value = -4604.937399999999
$ws.Range("N141").Value = -7710737.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3234.59
$ws.Range("I32").Value = 2565.4834
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2565.4834
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -2278.4834
$ws.Range("N32").Value = -10574
$ws.Range("H110").Value = 3077.7778
$ws.Range("I110").Value = 675
$ws.Range("K110").Value = 675
$ws.Range("M110").Value = 1370
$ws.Range("H122").Value = 2170.2424
$ws.Range("I122").Value = 1661.9131
$ws.Range("J122").Value = 3339.4
$ws.Range("K122").Value = 4985.7393
$ws.Range("L122").Value = 10018.2
$ws.Range("M122").Value = -2535.7393
$ws.Range("N122").Value = -14918.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2766.6584
$ws.Range("I134").Value = 2160.2
$ws.Range("J134").Value = 6304.3335
$ws.Range("K134").Value = 6480.599999999999
$ws.Range("L134").Value = 18913.0005
$ws.Range("M134").Value = -3945.599999999999
$ws.Range("N134").Value = -23983.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2746.6606
$ws.Range("I31").Value = 1563.5264
$ws.Range("J31").Value = 5244.3887
$ws.Range("K31").Value = 1563.5264
$ws.Range("L31").Value = 5244.3887
$ws.Range("M31").Value = -1268.5264
$ws.Range("N31").Value = -5834.3887
$ws.Range("H34").Value = 2746.6606
$ws.Range("I34").Value = 1563.5264
$ws.Range("J34").Value = 5244.3887
$ws.Range("K34").Value = 1563.5264
$ws.Range("L34").Value = 5244.3887
$ws.Range("M34").Value = -1361.5264
$ws.Range("N34").Value = -5648.3887
$ws.Range("H132").Value = 1393.2192
$ws.Range("I132").Value = 1028.0807
$ws.Range("J132").Value = 3451.2727
$ws.Range("K132").Value = 3084.2421
$ws.Range("L132").Value = 10353.8181
$ws.Range("M132").Value = -554.2420999999999
$ws.Range("N132").Value = -15413.8181
$ws.Range("H134").Value = 1498.3673
$ws.Range("I134").Value = 967.1515000000001
$ws.Range("J134").Value = 2594
$ws.Range("K134").Value = 2901.4545
$ws.Range("L134").Value = 7782
$ws.Range("M134").Value = -366.4545000000003
$ws.Range("N134").Value = -12852

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9279.9375
$ws.Range("I87").Value = 5197.6665
$ws.Range("J87").Value = 14528.571
$ws.Range("K87").Value = 15592.9995
$ws.Range("L87").Value = 43585.713
$ws.Range("M87").Value = -14344.9995
$ws.Range("N87").Value = -46081.713
$ws.Range("H90").Value = 9279.9375
$ws.Range("I90").Value = 5197.6665
$ws.Range("J90").Value = 14528.571
$ws.Range("K90").Value = 46778.9985
$ws.Range("L90").Value = 130757.139
$ws.Range("M90").Value = -40538.9985
$ws.Range("N90").Value = -143237.139
$ws.Range("H110").Value = 3330
$ws.Range("I110").Value = 965
$ws.Range("J110").Value = 3855.5557
$ws.Range("K110").Value = 2895
$ws.Range("L110").Value = 11566.6671
$ws.Range("M110").Value = 1195
$ws.Range("N110").Value = -19746.6671
$ws.Range("H118").Value = 1525.8334
$ws.Range("I118").Value = 475.7143
$ws.Range("J118").Value = 2996
$ws.Range("K118").Value = 1427.1429
$ws.Range("L118").Value = 8988
$ws.Range("M118").Value = -184.1428999999998
$ws.Range("N118").Value = -11474
$ws.Range("H120").Value = 19069.889
$ws.Range("I120").Value = 16515
$ws.Range("J120").Value = 19799.857
$ws.Range("K120").Value = 49545
$ws.Range("L120").Value = 59399.571
$ws.Range("M120").Value = -44707
$ws.Range("N120").Value = -69075.571
$ws.Range("H124").Value = 17371.428
$ws.Range("I124").Value = 1500
$ws.Range("J124").Value = 38533.332
$ws.Range("K124").Value = 4500
$ws.Range("L124").Value = 115599.996
$ws.Range("M124").Value = 410
$ws.Range("N124").Value = -125419.996
$ws.Range("H133").Value = 3783.7273
$ws.Range("J133").Value = 3380.25
$ws.Range("L133").Value = 10140.75
$ws.Range("N133").Value = -20260.75
$ws.Range("H134").Value = 2452
$ws.Range("I134").Value = 1595.3334
$ws.Range("J134").Value = 3994
$ws.Range("K134").Value = 4786.0002
$ws.Range("L134").Value = 11982
$ws.Range("M134").Value = 283.9997999999996
$ws.Range("N134").Value = -22122
$ws.Range("H140").Value = 9261275
$ws.Range("I140").Value = 12821472
$ws.Range("J140").Value = 4760
$ws.Range("K140").Value = 38464416
$ws.Range("L140").Value = 14280
$ws.Range("M140").Value = -38459236
$ws.Range("N140").Value = -24640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4530.207
$ws.Range("I70").Value = 4709.85
$ws.Range("J70").Value = 4131
$ws.Range("K70").Value = 4709.85
$ws.Range("L70").Value = 4131
$ws.Range("M70").Value = -4439.85
$ws.Range("N70").Value = -4671
$ws.Range("H73").Value = 4530.207
$ws.Range("I73").Value = 4709.85
$ws.Range("J73").Value = 4131
$ws.Range("K73").Value = 4709.85
$ws.Range("L73").Value = 4131
$ws.Range("M73").Value = -3773.85
$ws.Range("N73").Value = -6003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 14900
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 14900
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 14900
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -15240
$ws.Range("H19").Value = 6980
$ws.Range("I19").Value = 6980
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 6980
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -6810
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 125001624
$ws.Range("I22").Value = 166667090
$ws.Range("J22").Value = 5250
$ws.Range("K22").Value = 166667090
$ws.Range("L22").Value = 5250
$ws.Range("M22").Value = -166666795
$ws.Range("N22").Value = -5840
$ws.Range("H27").Value = 125001624
$ws.Range("I27").Value = 166667090
$ws.Range("J27").Value = 5250
$ws.Range("K27").Value = 166667090
$ws.Range("L27").Value = 5250
$ws.Range("M27").Value = -166666983
$ws.Range("N27").Value = -5464
$ws.Range("H46").Value = 1487.826
$ws.Range("I46").Value = 463.33334
$ws.Range("K46").Value = 463.33334
$ws.Range("M46").Value = -275.33334
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21322
$ws.Range("H132").Value = 1629.1528
$ws.Range("I132").Value = 1010.36206
$ws.Range("J132").Value = 4192.7144
$ws.Range("K132").Value = 3031.08618
$ws.Range("L132").Value = 12578.1432
$ws.Range("M132").Value = -501.0861800000002
$ws.Range("N132").Value = -17638.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14922.279
$ws.Range("I132").Value = 4083.6667
$ws.Range("J132").Value = 25268.227
$ws.Range("K132").Value = 12251.0001
$ws.Range("L132").Value = 75804.681
$ws.Range("M132").Value = -9721.000100000001
$ws.Range("N132").Value = -80864.681
$ws.Range("H136").Value = 1059.4131
$ws.Range("I136").Value = 492.59375
$ws.Range("J136").Value = 2355
$ws.Range("K136").Value = 1477.78125
$ws.Range("L136").Value = 7065
$ws.Range("M136").Value = 1072.21875
$ws.Range("N136").Value = -12165
